$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 36.333332
$ws.Range("I11").Value = 36.333332
$ws.Range("K11").Value = 36.333332
$ws.Range("M11").Value = 103.666668
$ws.Range("H17").Value = 2349.75
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H39").Value = 314.33334
$ws.Range("I39").Value = 224.83333
$ws.Range("K39").Value = 674.49999
$ws.Range("M39").Value = -378.49999
$ws.Range("H132").Value = 1235.3334
$ws.Range("I132").Value = 1235.3334
$ws.Range("K132").Value = 3706.0002
$ws.Range("M132").Value = -1176.0002
$ws.Range("H135").Value = 1534.4667
$ws.Range("I135").Value = 1534.4667
$ws.Range("K135").Value = 13810.2003
$ws.Range("M135").Value = -11275.2003
$ws.Range("H137").Value = 605780.5
$ws.Range("I137").Value = 1378.6666
$ws.Range("K137").Value = 4135.9998
$ws.Range("M137").Value = -1585.9998
$ws.Range("H138").Value = 1616.9302
$ws.Range("I138").Value = 1147.6428
$ws.Range("K138").Value = 3442.9284
$ws.Range("M138").Value = 1697.0716
$ws.Range("H141").Value = 3268.182
$ws.Range("I141").Value = 3268.182
$ws.Range("K141").Value = 9804.545999999998
$ws.Range("M141").Value = -4624.545999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1595.579
$ws.Range("I2").Value = 1497.6428
$ws.Range("K2").Value = 1497.6428
$ws.Range("M2").Value = -1384.6428
$ws.Range("H61").Value = 2785.7144
$ws.Range("I61").Value = 2100
$ws.Range("J61").Value = 3700
$ws.Range("K61").Value = 2100
$ws.Range("L61").Value = 3700
$ws.Range("M61").Value = -1888
$ws.Range("N61").Value = -4124
$ws.Range("H74").Value = 1858.5625
$ws.Range("I74").Value = 1356.9445
$ws.Range("K74").Value = 1356.9445
$ws.Range("M74").Value = -482.9445000000001
$ws.Range("H76").Value = 124998
$ws.Range("J76").Value = 124998
$ws.Range("L76").Value = 124998
$ws.Range("N76").Value = -125674
$ws.Range("H77").Value = 1858.5625
$ws.Range("I77").Value = 1356.9445
$ws.Range("K77").Value = 6784.7225
$ws.Range("M77").Value = -2416.7225
$ws.Range("H79").Value = 124998
$ws.Range("J79").Value = 124998
$ws.Range("L79").Value = 124998
$ws.Range("N79").Value = -127338
$ws.Range("H116").Value = 1595.579
$ws.Range("I116").Value = 1497.6428
$ws.Range("K116").Value = 1497.6428
$ws.Range("M116").Value = 796.3571999999999
$ws.Range("H122").Value = 2833.6875
$ws.Range("I122").Value = 2553.25
$ws.Range("K122").Value = 7659.75
$ws.Range("M122").Value = -5209.75
$ws.Range("H132").Value = 1778.2632
$ws.Range("I132").Value = 1639.9354
$ws.Range("K132").Value = 4919.8062
$ws.Range("M132").Value = -2389.8062
$ws.Range("H136").Value = 2785.7144
$ws.Range("I136").Value = 2100
$ws.Range("J136").Value = 3700
$ws.Range("K136").Value = 6300
$ws.Range("L136").Value = 11100
$ws.Range("M136").Value = -3750
$ws.Range("N136").Value = -16200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1595.579
$ws.Range("I3").Value = 1497.6428
$ws.Range("K3").Value = 1497.6428
$ws.Range("M3").Value = -1383.6428
$ws.Range("H134").Value = 2188.3333
$ws.Range("I134").Value = 1599.4286
$ws.Range("J134").Value = 4249.5
$ws.Range("K134").Value = 4798.2858
$ws.Range("L134").Value = 12748.5
$ws.Range("M134").Value = -2263.2858
$ws.Range("N134").Value = -17818.5
$ws.Range("H140").Value = 101909.086
$ws.Range("J140").Value = 65719
$ws.Range("L140").Value = 65719
$ws.Range("N140").Value = -76079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 35889.816
$ws.Range("J9").Value = 35889.816
$ws.Range("L9").Value = 35889.816
$ws.Range("N9").Value = -36225.816
$ws.Range("H74").Value = 44000
$ws.Range("J74").Value = 44000
$ws.Range("L74").Value = 44000
$ws.Range("N74").Value = -45748
$ws.Range("H77").Value = 44000
$ws.Range("J77").Value = 44000
$ws.Range("L77").Value = 132000
$ws.Range("N77").Value = -140736
$ws.Range("H122").Value = 3168.8572
$ws.Range("I122").Value = 2118.2222
$ws.Range("K122").Value = 6354.6666
$ws.Range("M122").Value = -3904.6666
$ws.Range("H132").Value = 4002500.2
$ws.Range("I132").Value = 4547830
$ws.Range("J132").Value = 3416.6667
$ws.Range("K132").Value = 13643490
$ws.Range("L132").Value = 10250.0001
$ws.Range("M132").Value = -13640960
$ws.Range("N132").Value = -15310.0001
$ws.Range("H134").Value = 1840.5454
$ws.Range("I134").Value = 1147.4286
$ws.Range("J134").Value = 3053.5
$ws.Range("K134").Value = 3442.2858
$ws.Range("L134").Value = 9160.5
$ws.Range("M134").Value = -907.2857999999997
$ws.Range("N134").Value = -14230.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1323.579
$ws.Range("I5").Value = 679.9167
$ws.Range("J5").Value = 2427
$ws.Range("K5").Value = 2039.7501
$ws.Range("L5").Value = 7281
$ws.Range("M5").Value = -1927.7501
$ws.Range("N5").Value = -7505
$ws.Range("H44").Value = 1050
$ws.Range("I44").Value = 100
$ws.Range("J44").Value = 2000
$ws.Range("K44").Value = 300
$ws.Range("L44").Value = 6000
$ws.Range("M44").Value = 98
$ws.Range("N44").Value = -6796
$ws.Range("H50").Value = 481.2857
$ws.Range("I50").Value = 517.6667
$ws.Range("K50").Value = 1553.0001
$ws.Range("M50").Value = -1072.0001
$ws.Range("H53").Value = 481.2857
$ws.Range("I53").Value = 517.6667
$ws.Range("K53").Value = 1553.0001
$ws.Range("M53").Value = -1072.0001
$ws.Range("H68").Value = 2200
$ws.Range("J68").Value = 2200
$ws.Range("L68").Value = 6600
$ws.Range("N68").Value = -8222
$ws.Range("H71").Value = 2200
$ws.Range("J71").Value = 2200
$ws.Range("L71").Value = 19800
$ws.Range("N71").Value = -27912
$ws.Range("H135").Value = 1323.579
$ws.Range("I135").Value = 679.9167
$ws.Range("J135").Value = 2427
$ws.Range("K135").Value = 6119.2503
$ws.Range("L135").Value = 21843
$ws.Range("M135").Value = -3584.2503
$ws.Range("N135").Value = -26913

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 142859860
$ws.Range("I80").Value = 500001440
$ws.Range("J80").Value = 3208
$ws.Range("K80").Value = 500001440
$ws.Range("L80").Value = 3208
$ws.Range("M80").Value = -500000442
$ws.Range("N80").Value = -5204
$ws.Range("H83").Value = 142859860
$ws.Range("I83").Value = 500001440
$ws.Range("J83").Value = 3208
$ws.Range("K83").Value = 2500007200
$ws.Range("L83").Value = 16040
$ws.Range("M83").Value = -2500002208
$ws.Range("N83").Value = -26024
$ws.Range("H102").Value = 1524.8334
$ws.Range("I102").Value = 1481.6364
$ws.Range("K102").Value = 1481.6364
$ws.Range("M102").Value = 140.3635999999999
$ws.Range("H132").Value = 3724.851
$ws.Range("I132").Value = 3073.394
$ws.Range("J132").Value = 5260.4287
$ws.Range("K132").Value = 9220.181999999999
$ws.Range("L132").Value = 15781.2861
$ws.Range("M132").Value = -6690.181999999999
$ws.Range("N132").Value = -20841.2861
$ws.Range("H134").Value = 31841.75
$ws.Range("J134").Value = 31841.75
$ws.Range("L134").Value = 95525.25
$ws.Range("N134").Value = -100595.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5559256.5
$ws.Range("I40").Value = 4126.25
$ws.Range("K40").Value = 4126.25
$ws.Range("M40").Value = -3990.25
$ws.Range("H122").Value = 15412486
$ws.Range("I122").Value = 38234
$ws.Range("K122").Value = 114702
$ws.Range("M122").Value = -112252
$ws.Range("H132").Value = 15149.546
$ws.Range("I132").Value = 24658.334
$ws.Range("J132").Value = 3739
$ws.Range("K132").Value = 73975.00199999999
$ws.Range("L132").Value = 11217
$ws.Range("M132").Value = -71445.00199999999
$ws.Range("N132").Value = -16277
$ws.Range("H136").Value = 4734.0977
$ws.Range("I136").Value = 3905.4075
$ws.Range("J136").Value = 6332.2856
$ws.Range("K136").Value = 11716.2225
$ws.Range("L136").Value = 18996.8568
$ws.Range("M136").Value = -9166.2225
$ws.Range("N136").Value = -24096.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1360224
$ws.Range("I132").Value = 1445.069
$ws.Range("J132").Value = 14495087
$ws.Range("K132").Value = 4335.207
$ws.Range("L132").Value = 43485261
$ws.Range("M132").Value = -1805.207
$ws.Range("N132").Value = -43490321
$ws.Range("H136").Value = 3497.12
$ws.Range("I136").Value = 1771.45
$ws.Range("J136").Value = 10399.8
$ws.Range("K136").Value = 5314.35
$ws.Range("L136").Value = 31199.4
$ws.Range("M136").Value = -2764.35
$ws.Range("N136").Value = -36299.39999999999
